$ErrorActionPreference = "Stop"
$d = $word.ActiveDocument

function Get-PkgWrapped($innerBodyXml) {
    return '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $innerBodyXml + '</w:document></pkg:xmlData></pkg:part></pkg:package>'
}

function Find-ParagraphContaining($searchText) {
    $r = $d.Content
    $found = $r.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Text not found: $searchText"
    }
    $fStart = $r.Start
    $fEnd = $r.End
    $allParas = $d.Paragraphs
    for ($i = 1; $i -le $allParas.Count; $i++) {
        $pp = $allParas.Item($i)
        if ($pp.Range.Start -le $fStart -and $pp.Range.End -ge $fEnd) {
            return $pp
        }
    }
    throw "Paragraph not found for: $searchText"
}

function Replace-ParagraphXml($paragraph, $newParagraphXml) {
    $rng = $paragraph.Range
    $rng.InsertXML( (Get-PkgWrapped ("<w:body>" + $newParagraphXml + "</w:body>")) )
}

# 1) "LoginScreen" cell: drop the spell-check proofErr wrapper around the run (no visible text change)
$pLogin = Find-ParagraphContaining "LoginScreen"
Replace-ParagraphXml $pLogin '<w:p><w:pPr><w:cnfStyle w:val="000000100000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="1" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/></w:pPr><w:r><w:t>LoginScreen</w:t></w:r></w:p>'

# 2) "User kann sich mit vergebenen Logindaten einloggen" cell: merge the 3 runs (with proofErr around "Logindaten") into one run
$pLogindaten = Find-ParagraphContaining "Logindaten"
Replace-ParagraphXml $pLogindaten '<w:p><w:pPr><w:cnfStyle w:val="000000100000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="1" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/></w:pPr><w:r><w:t>User kann sich mit vergebenen Logindaten einloggen</w:t></w:r></w:p>'

# 3) Remove the old "_GoBack" bookmark from the FA6/Verlegung row's empty "Erfolgreich?" cell
$bm = $d.Bookmarks.Item("_GoBack")
$bmStart = $bm.Start
$allParas = $d.Paragraphs
$pBookmark = $null
for ($i = 1; $i -le $allParas.Count; $i++) {
    $pp = $allParas.Item($i)
    if ($pp.Range.Start -eq $bmStart -and $pp.Range.End -gt $pp.Range.Start) {
        $pBookmark = $pp
        break
    }
}
if ($null -eq $pBookmark) {
    throw "Could not locate paragraph hosting the _GoBack bookmark"
}
Replace-ParagraphXml $pBookmark '<w:p><w:pPr><w:cnfStyle w:val="000000100000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="1" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/></w:pPr></w:p>'

# 4) Append the new "FA8 - Entlassung" heading + table (with the new, re-homed "_GoBack" bookmark)
#    after the last table, between the two trailing empty paragraphs.
#    (Inserting XML whose last block is a <w:tbl> right at the very last paragraph of the
#    document swallows that trailing paragraph instead of pushing it down, so first grow
#    the tail by one empty paragraph via the object model, then InsertXML into the new slot -
#    that keeps both original trailing empty paragraphs intact.)
$n = $d.Paragraphs.Count
$secondLast = $d.Paragraphs.Item($n - 1)
$growRng = $secondLast.Range
$growRng.Collapse(0)
$growRng.InsertParagraphAfter()

$headingPara = $d.Paragraphs.Item($n)
$insertionPoint = $headingPara.Range
$insertionPoint.Collapse(1)
$fa8Xml = @'
<w:p>
  <w:pPr>
    <w:pStyle w:val="berschrift1"/>
  </w:pPr>
  <w:r>
    <w:t>FA</w:t>
  </w:r>
  <w:r>
    <w:t>8</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> – </w:t>
  </w:r>
  <w:r>
    <w:t>Entlassung</w:t>
  </w:r>
</w:p>
<w:tbl>
  <w:tblPr>
    <w:tblStyle w:val="HellesRaster-Akzent3"/>
    <w:tblW w:w="11058" w:type="dxa"/>
    <w:tblInd w:w="-885" w:type="dxa"/>
    <w:tblLayout w:type="fixed"/>
    <w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/>
  </w:tblPr>
  <w:tblGrid>
    <w:gridCol w:w="993"/>
    <w:gridCol w:w="2127"/>
    <w:gridCol w:w="4961"/>
    <w:gridCol w:w="1559"/>
    <w:gridCol w:w="1418"/>
  </w:tblGrid>
  <w:tr>
    <w:trPr>
      <w:cnfStyle w:val="100000000000" w:firstRow="1" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/>
    </w:trPr>
    <w:tc>
      <w:tcPr>
        <w:cnfStyle w:val="001000000000" w:firstRow="0" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/>
        <w:tcW w:w="993" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:t>Nr.</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="2127" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:cnfStyle w:val="100000000000" w:firstRow="1" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/>
        </w:pPr>
        <w:r>
          <w:t>Testbereich</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="4961" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:cnfStyle w:val="100000000000" w:firstRow="1" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/>
        </w:pPr>
        <w:r>
          <w:t>Erwartetes Ergebnis</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="1559" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:cnfStyle w:val="100000000000" w:firstRow="1" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/>
        </w:pPr>
        <w:r>
          <w:t>Erfolgreich?</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="1418" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:cnfStyle w:val="100000000000" w:firstRow="1" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/>
        </w:pPr>
        <w:r>
          <w:t>Testdatum</w:t>
        </w:r>
      </w:p>
    </w:tc>
  </w:tr>
  <w:tr>
    <w:trPr>
      <w:cnfStyle w:val="000000100000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="1" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/>
    </w:trPr>
    <w:tc>
      <w:tcPr>
        <w:cnfStyle w:val="001000000000" w:firstRow="0" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/>
        <w:tcW w:w="993" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:t>FA8.1</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="2127" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:cnfStyle w:val="000000100000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="1" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/>
        </w:pPr>
        <w:r>
          <w:t>Patientenverwaltung</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="4961" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:cnfStyle w:val="000000100000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="1" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/>
        </w:pPr>
        <w:r>
          <w:t>Patient kann entlassen werden und wird dann nicht mehr vom System gefunden</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="1559" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:cnfStyle w:val="000000100000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="1" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/>
        </w:pPr>
        <w:r>
          <w:t>ja</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="1418" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:cnfStyle w:val="000000100000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="1" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/>
        </w:pPr>
        <w:r>
          <w:t>04.06.2019</w:t>
        </w:r>
      </w:p>
    </w:tc>
  </w:tr>
  <w:tr>
    <w:trPr>
      <w:cnfStyle w:val="000000010000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="1" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/>
    </w:trPr>
    <w:tc>
      <w:tcPr>
        <w:cnfStyle w:val="001000000000" w:firstRow="0" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/>
        <w:tcW w:w="993" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:t>FA</w:t>
        </w:r>
        <w:r>
          <w:t>8</w:t>
        </w:r>
        <w:r>
          <w:t>.2</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="2127" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:cnfStyle w:val="000000010000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="1" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/>
        </w:pPr>
        <w:r>
          <w:t>Datenbank</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="4961" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:cnfStyle w:val="000000010000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="1" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/>
        </w:pPr>
        <w:r>
          <w:t>Bett wurde wieder frei gegeben</w:t>
        </w:r>
        <w:bookmarkStart w:id="0" w:name="_GoBack"/>
        <w:bookmarkEnd w:id="0"/>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="1559" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:cnfStyle w:val="000000010000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="1" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/>
        </w:pPr>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="1418" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:cnfStyle w:val="000000010000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="1" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/>
        </w:pPr>
      </w:p>
    </w:tc>
  </w:tr>
</w:tbl>

'@
$insertionPoint.InsertXML( (Get-PkgWrapped $fa8Xml) )

Write-Host "All edits applied."
